{"js": "// Append the 13/11 journal entry block to the end of the document, then\n// move the Word \"last edit\" (_GoBack) bookmark from the old final\n// paragraph to the end of the new final content paragraph, and finish\n// with one trailing blank paragraph - matching the author's commit.\n\nconst body = context.document.body;\n\n// The existing \"_GoBack\" bookmark sits at the end of the previously-last\n// paragraph (\"...revisited in later sprints.\"). As the journal entry grows\n// Word keeps this bookmark pinned to the newest edit, so drop it here; it\n// gets re-created further down once the new text is in place.\ncontext.document.deleteBookmark(\"_GoBack\");\n\n// Five blank lines, then the date heading - all inherit the surrounding\n// \"firstLine 720 twip\" indent automatically from the preceding paragraph.\nbody.insertParagraph(\"\", \"End\");\nbody.insertParagraph(\"\", \"End\");\nbody.insertParagraph(\"\", \"End\");\nbody.insertParagraph(\"\", \"End\");\nbody.insertParagraph(\"\", \"End\");\nbody.insertParagraph(\"13/11\", \"End\");\n\nbody.insertParagraph(\n  \"Refactored query lists remove significant amounts of redundant and repeated code\",\n  \"End\"\n);\n\nbody.insertParagraph(\n  \"Rearranged data in query invite process requirement of values to validated authorative changes moved to member / pending invite objects/\",\n  \"End\"\n);\n\nbody.insertParagraph(\"\", \"End\");\n\nbody.insertParagraph(\"Significant issues with array modification \", \"End\");\n\nbody.insertParagraph(\n  \"I have gone significantly over my estimate for this period of time due to unforeseen refactoring needed, I should have done UML diagrams earlier and made more satalight tests.\",\n  \"End\"\n);\nawait context.sync();\n\n// Re-fetch the paragraph collection so the last item is a \"live\" proxy -\n// re-anchor \"_GoBack\" right after the final run of text (but before the\n// paragraph mark) in this paragraph, which is where Word leaves it after\n// the last keystroke of the session.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\nconst lastContentParagraph = paragraphs.items[paragraphs.items.length - 1];\nlastContentParagraph.getRange(\"End\").insertBookmark(\"_GoBack\");\nawait context.sync();\n\n// The author left one more empty line after their last sentence.\nbody.insertParagraph(\"\", \"End\");\n\nawait context.sync();\n", "ps1": "# Append the 13/11 journal entry block to the end of the document, then\n# move the Word \"last edit\" (_GoBack) bookmark from the old final\n# paragraph to the end of the new final content paragraph, and finish\n# with one trailing blank paragraph - matching the author's commit.\n\n$d = $word.ActiveDocument\n\n# The existing \"_GoBack\" bookmark sits at the end of the previously-last\n# paragraph (\"...revisited in later sprints.\"). As the journal entry grows\n# Word keeps this bookmark pinned to the newest edit, so drop it here; it\n# gets re-created further down once the new text is in place.\n$d.Bookmarks.Item(\"_GoBack\").Delete()\n\n# Five blank lines, then the date heading - all inherit the surrounding\n# \"firstLine 720 twip\" indent automatically from the preceding paragraph.\n$d.Paragraphs.Last.Range.InsertParagraphAfter()\n$d.Paragraphs.Last.Range.InsertParagraphAfter()\n$d.Paragraphs.Last.Range.InsertParagraphAfter()\n$d.Paragraphs.Last.Range.InsertParagraphAfter()\n$d.Paragraphs.Last.Range.InsertParagraphAfter()\n\n$d.Paragraphs.Last.Range.InsertParagraphAfter()\n$d.Paragraphs.Last.Range.Text = \"13/11\"\n\n$d.Paragraphs.Last.Range.InsertParagraphAfter()\n$d.Paragraphs.Last.Range.Text = \"Refactored query lists remove significant amounts of redundant and repeated code\"\n\n$d.Paragraphs.Last.Range.InsertParagraphAfter()\n$d.Paragraphs.Last.Range.Text = \"Rearranged data in query invite process requirement of values to validated authorative changes moved to member / pending invite objects/\"\n\n$d.Paragraphs.Last.Range.InsertParagraphAfter()\n\n$d.Paragraphs.Last.Range.InsertParagraphAfter()\n$d.Paragraphs.Last.Range.Text = \"Significant issues with array modification \"\n\n$d.Paragraphs.Last.Range.InsertParagraphAfter()\n$d.Paragraphs.Last.Range.Text = \"I have gone significantly over my estimate for this period of time due to unforeseen refactoring needed, I should have done UML diagrams earlier and made more satalight tests.\"\n\n# Re-anchor \"_GoBack\" right after the final run of text (but before the\n# paragraph mark) in this paragraph, which is where Word leaves it after\n# the last keystroke of the session. Shrinking the paragraph range by one\n# character (the trailing paragraph mark) before adding the bookmark keeps\n# it collapsed at that point instead of spanning the whole paragraph.\n$lastContentParagraph = $d.Paragraphs.Last\n$bookmarkRange = $lastContentParagraph.Range\n$bookmarkRange.MoveEnd(1, -1)\n$d.Bookmarks.Add(\"_GoBack\", $bookmarkRange)\n\n# The author left one more empty line after their last sentence.\n$d.Paragraphs.Last.Range.InsertParagraphAfter()\n"}
